$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '28.031.61'
$ws.Range('E2').Value = '  -0.71%  '
$ws.Range('D3').Value = '1.751.01'
$ws.Range('E3').Value = '  -1.90%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.9965'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.61%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '333.58'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.28%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9975'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3879'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +1.26%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3390'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.61%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '45.32'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.28%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.112'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.98%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07194'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -3.29%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.9933'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.65%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '22.10'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -5.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.115'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -5.09%  '
$ws.Range('D15').Value = '1.741.65'
$ws.Range('E15').Value = '  -2.44%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.986'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.19%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001054'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -2.70%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.06598'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.30%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '79.76'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -3.41%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.9977'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.15%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '16.85'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -3.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.188'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -4.29%  '
$ws.Range('D23').Value = '27.946.53'
$ws.Range('E23').Value = '  -1.03%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.55'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -4.67%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.372'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '153.72'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -1.34%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '19.87'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -4.90%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '2.304'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -5.01%  '
$ws.Range('D29').Value = '1.939.75'
$ws.Range('E29').Value = '  -2.34%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.275'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -11.51%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '128.25'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.59%  '
$ws.Range('E32').Value = '  +2.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.793'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -5.95%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08666'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.97%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '12.04'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -6.35%  '
$ws.Range('B36').Value = 'VeChain'
$ws.Range('C36').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.02277'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.63%  '
$ws.Range('B37').Value = 'InternetComputer(DFINITY)'
$ws.Range('C37').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '5.119'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.33%  '
$ws.Range('B38').Value = 'Hedera'
$ws.Range('C38').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.06110'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -4.02%  '
$ws.Range('B39').Value = 'WEMIXTOKEN'
$ws.Range('C39').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.508'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +0.48%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.6448'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -6.23%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.2097'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -4.15%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.197'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -3.74%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.0000'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '7.923'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -4.94%  '
$ws.Range('E45').Value = '  -5.49%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.808'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -1.56%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5954'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -5.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '126.01'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -5.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.980'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06956'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -6.81%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.146'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -5.37%  '
